# Commit: "Fruta / hortaliza, semanal"
#
# A new weekly price record for Apio (Macroferia Regional de Talca) is
# inserted as row 127, pushing the existing rows 127..242 down to 128..243
# (the last former row, 242, becomes row 243). Columns A, B, C, E, F, G, H,
# I and R are identical on every data row of this sheet, so only D (Fecha),
# J (Volumen), K/L/M (Precio min/max/promedio), N (Unidad), O (Origen),
# P (Precio $/Kg) and Q (Kg o Unidades) need explicit values for the new
# record; the rest are copied straight from the template row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 127:242 down to 128:243, carrying their values/styles with them.
$ws.Rows.Item(127).Insert()

# Populate the newly inserted row 127 with the new weekly record.
$ws.Cells.Item(127, 1).Value = 5
$ws.Cells.Item(127, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(127, 3).Value = "Maule"
$ws.Cells.Item(127, 4).Value = 44827
$ws.Cells.Item(127, 5).Value = 7
$ws.Cells.Item(127, 6).Value = 100112017
$ws.Cells.Item(127, 7).Value = "Apio"
$ws.Cells.Item(127, 8).Value = "Americana (o)"
$ws.Cells.Item(127, 9).Value = "Primera"
$ws.Cells.Item(127, 10).Value = 700
$ws.Cells.Item(127, 11).Value = 9000
$ws.Cells.Item(127, 12).Value = 9000
$ws.Cells.Item(127, 13).Value = 9000
$ws.Cells.Item(127, 14).Value = "`$/docena de matas"
$ws.Cells.Item(127, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(127, 16).Value = 1500
$ws.Cells.Item(127, 17).Value = 6
$ws.Cells.Item(127, 18).Value = "Hortaliza"
